$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values are entered in the same order the original author typed them so the
# shared-strings table is built in a matching sequence.
$ws.Range("A10").Value = "First.csv"
$ws.Range("A9").Value = "FILES"
$ws.Range("B9").Value = "INFO"
$ws.Range("B10").Value = "Drive behind Switch"

$ws.Range("A11").Value = "Switch_Close.csv"
$ws.Range("A12").Value = "Switch_Far.csv"
$ws.Range("A13").Value = "Scale_Close.csv"
$ws.Range("A14").Value = "Scale_Far.csv"

$ws.Range("B13").Value = "Score in scale platform directly behind robot"
$ws.Range("B14").Value = "Score in scale platform across from robot"
$ws.Range("B12").Value = "Score in far switch platform & collect block"
$ws.Range("B11").Value = "Score in close switch platform & collect block"

$ws.Range("A9:B9").Font.Bold = $true

$ws.Range("E14").Select()
